$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and G (Hora) contain numeric-looking text that must
# stay stored as text (matching the original inlineStr cell type), so we
# force a text number format before assigning the value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.30"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "20"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.76"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "20"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.196"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "20"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06099"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "20"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.514"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "20"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.725"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "20"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.358"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "20"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.7983"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "20"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "20"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08087"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "20"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03342"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "20"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03101"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "20"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09274"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "20"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.890"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "20"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001695"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "20"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "20"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006205"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "20"
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.001100"
$ws.Range("E19").Value = "18BitKanKAN"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "20"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.003398"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "20"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001500"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "20"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.691"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "20"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.260"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "20"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01332"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "20"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3362"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "20"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1226"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "20"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "20"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "20"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "20"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "20"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "20"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "20"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "20"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "20"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "20"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "20"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "20"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "20"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "20"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04592"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "20"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007087"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "20"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003899"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "20"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1122"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "20"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01023"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "20"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.002970"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "20"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006005"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "20"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "20"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "20"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05948"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "20"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "20"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "20"
